# Update the "想去人数" (F column) values for the 展览 and 全部类型 sheets,
# matching the regenerated data snapshot.

$wb = $excel.ActiveWorkbook

# Map of row -> new value for column F
$updates = @{
    2  = 1164
    3  = 97
    4  = 1547
    5  = 592
    7  = 11240
    8  = 8
    10 = 220
    12 = 1079
    13 = 774
    14 = 12286
    15 = 12903
    17 = 133
    20 = 80
    22 = 69
}

# Sheets "展览" (index 1) and "全部类型" (index 4) both contain this table
# and both received identical updates in this commit.
$sheetIndexes = @(1, 4)

foreach ($sheetIndex in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
